$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 125.68
$ws.Range("I15").Value = 125.68
$ws.Range("K15").Value = 377.04
$ws.Range("M15").Value = -208.04

$ws.Range("H132").Value = 16460806
$ws.Range("I132").Value = 19309126
$ws.Range("J132").Value = 3842.7778
$ws.Range("K132").Value = 57927378
$ws.Range("L132").Value = 11528.3334
$ws.Range("M132").Value = -57924848
$ws.Range("N132").Value = -16588.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1092.5
$ws.Range("I2").Value = 1085.8
$ws.Range("J2").Value = 1103.6666
$ws.Range("K2").Value = 1085.8
$ws.Range("L2").Value = 1103.6666
$ws.Range("M2").Value = -972.8
$ws.Range("N2").Value = -1329.6666

$ws.Range("H28").Value = 15500
$ws.Range("I28").Value = 15500
$ws.Range("K28").Value = 15500
$ws.Range("M28").Value = -15308

$ws.Range("H99").Value = 15500
$ws.Range("I99").Value = 15500
$ws.Range("K99").Value = 15500
$ws.Range("M99").Value = -12505

$ws.Range("H116").Value = 1092.5
$ws.Range("I116").Value = 1085.8
$ws.Range("J116").Value = 1103.6666
$ws.Range("K116").Value = 1085.8
$ws.Range("L116").Value = 1103.6666
$ws.Range("M116").Value = 1208.2
$ws.Range("N116").Value = -5691.6666

$ws.Range("H132").Value = 2164.5117
$ws.Range("I132").Value = 1050.0625
$ws.Range("J132").Value = 5406.5454
$ws.Range("K132").Value = 3150.1875
$ws.Range("L132").Value = 16219.6362
$ws.Range("M132").Value = -620.1875
$ws.Range("N132").Value = -21279.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1092.5
$ws.Range("I3").Value = 1085.8
$ws.Range("J3").Value = 1103.6666
$ws.Range("K3").Value = 1085.8
$ws.Range("L3").Value = 1103.6666
$ws.Range("M3").Value = -971.8
$ws.Range("N3").Value = -1331.6666

$ws.Range("H26").Value = 20471
$ws.Range("I26").Value = 20471
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 20471
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -20179
$ws.Range("N26").ClearContents()

$ws.Range("H96").Value = 27809.334
$ws.Range("I96").Value = 24214
$ws.Range("J96").Value = 35000
$ws.Range("K96").Value = 24214
$ws.Range("L96").Value = 35000
$ws.Range("M96").Value = -21468
$ws.Range("N96").Value = -40492

$ws.Range("H105").Value = 2449
$ws.Range("I105").Value = 2205
$ws.Range("K105").Value = 2205
$ws.Range("M105").Value = -458

$ws.Range("H134").Value = 2971.1333
$ws.Range("I134").Value = 1583.0571
$ws.Range("J134").Value = 7829.4
$ws.Range("K134").Value = 4749.1713
$ws.Range("L134").Value = 23488.2
$ws.Range("M134").Value = -2214.1713
$ws.Range("N134").Value = -28558.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2902.5957
$ws.Range("I31").Value = 1350.6765
$ws.Range("J31").Value = 6961.4614
$ws.Range("K31").Value = 1350.6765
$ws.Range("L31").Value = 6961.4614
$ws.Range("M31").Value = -1055.6765
$ws.Range("N31").Value = -7551.4614

$ws.Range("H34").Value = 2902.5957
$ws.Range("I34").Value = 1350.6765
$ws.Range("J34").Value = 6961.4614
$ws.Range("K34").Value = 1350.6765
$ws.Range("L34").Value = 6961.4614
$ws.Range("M34").Value = -1148.6765
$ws.Range("N34").Value = -7365.4614

$ws.Range("H36").Value = 21881
$ws.Range("I36").Value = 10024
$ws.Range("K36").Value = 10024
$ws.Range("M36").Value = -9636

$ws.Range("H40").Value = 21881
$ws.Range("I40").Value = 10024
$ws.Range("K40").Value = 10024
$ws.Range("M40").Value = -9864

$ws.Range("H58").Value = 1960.9
$ws.Range("I58").Value = 1661.9375
$ws.Range("J58").Value = 5149.8335
$ws.Range("K58").Value = 1661.9375
$ws.Range("L58").Value = 5149.8335
$ws.Range("M58").Value = -1458.9375
$ws.Range("N58").Value = -5555.8335

$ws.Range("H134").Value = 4911.0645
$ws.Range("I134").Value = 4970.885
$ws.Range("J134").Value = 4600
$ws.Range("K134").Value = 14912.655
$ws.Range("L134").Value = 13800
$ws.Range("M134").Value = -12377.655
$ws.Range("N134").Value = -18870

$ws.Range("H136").Value = 1960.9
$ws.Range("I136").Value = 1661.9375
$ws.Range("J136").Value = 5149.8335
$ws.Range("K136").Value = 4985.8125
$ws.Range("L136").Value = 15449.5005
$ws.Range("M136").Value = -2435.8125
$ws.Range("N136").Value = -20549.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2196214.5
$ws.Range("J4").Value = 6098.25
$ws.Range("L4").Value = 18294.75
$ws.Range("N4").Value = -18518.75

$ws.Range("H39").Value = 8201.214
$ws.Range("J39").Value = 8616.615
$ws.Range("L39").Value = 25849.845
$ws.Range("N39").Value = -26437.845

$ws.Range("H58").Value = 2168.3333
$ws.Range("I58").Value = 2052.5
$ws.Range("J58").Value = 2400
$ws.Range("K58").Value = 6157.5
$ws.Range("L58").Value = 7200
$ws.Range("M58").Value = -6029.5
$ws.Range("N58").Value = -7456

$ws.Range("H76").Value = 3730.4285
$ws.Range("I76").Value = 2006.5
$ws.Range("J76").Value = 4420
$ws.Range("K76").Value = 6019.5
$ws.Range("L76").Value = 13260
$ws.Range("M76").Value = -5636.5
$ws.Range("N76").Value = -14026

$ws.Range("H79").Value = 3730.4285
$ws.Range("I79").Value = 2006.5
$ws.Range("J79").Value = 4420
$ws.Range("K79").Value = 6019.5
$ws.Range("L79").Value = 13260
$ws.Range("M79").Value = -4693.5
$ws.Range("N79").Value = -15912

$ws.Range("H112").Value = 4177.25
$ws.Range("I112").Value = 1913.5
$ws.Range("J112").Value = 4630
$ws.Range("K112").Value = 5740.5
$ws.Range("L112").Value = 13890
$ws.Range("M112").Value = -4632.5
$ws.Range("N112").Value = -16106

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 28857.143
$ws.Range("J4").Value = 28857.143
$ws.Range("L4").Value = 28857.143
$ws.Range("N4").Value = -29081.143

$ws.Range("H70").Value = 6499.154
$ws.Range("I70").Value = 5829.9443
$ws.Range("K70").Value = 5829.9443
$ws.Range("M70").Value = -5559.9443

$ws.Range("H73").Value = 6499.154
$ws.Range("I73").Value = 5829.9443
$ws.Range("K73").Value = 5829.9443
$ws.Range("M73").Value = -4893.9443

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1558.7428
$ws.Range("I22").Value = 1206
$ws.Range("K22").Value = 1206
$ws.Range("M22").Value = -911

$ws.Range("H27").Value = 1558.7428
$ws.Range("I27").Value = 1206
$ws.Range("K27").Value = 1206
$ws.Range("M27").Value = -1099

$ws.Range("H94").Value = 33750
$ws.Range("J94").Value = 33750
$ws.Range("L94").Value = 33750
$ws.Range("N94").Value = -35102

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 35600
$ws.Range("I42").Value = 2400
$ws.Range("J42").Value = 46666.668
$ws.Range("K42").Value = 2400
$ws.Range("L42").Value = 46666.668
$ws.Range("M42").Value = -2022
$ws.Range("N42").Value = -47422.668

$ws.Range("H92").Value = 32000
$ws.Range("J92").Value = 32000
$ws.Range("L92").Value = 32000
$ws.Range("N92").Value = -36992

$ws.Range("H93").Value = 32500
$ws.Range("J93").Value = 32500
$ws.Range("L93").Value = 32500
$ws.Range("N93").Value = -37492

$ws.Range("H103").Value = 35133.332
$ws.Range("J103").Value = 35133.332
$ws.Range("L103").Value = 35133.332
$ws.Range("N103").Value = -37477.332

$ws.Range("H104").Value = 27333.334
$ws.Range("J104").Value = 27333.334
$ws.Range("L104").Value = 27333.334
$ws.Range("N104").Value = -34321.334

$ws.Range("H106").Value = 18282.428
$ws.Range("J106").Value = 18282.428
$ws.Range("L106").Value = 18282.428
$ws.Range("N106").Value = -20806.428

$ws.Range("H118").Value = 28890
$ws.Range("J118").Value = 28890
$ws.Range("L118").Value = 28890
$ws.Range("N118").Value = -32204
